$wb = $excel.ActiveWorkbook

# "L3" sheet holds the groups list; it previously stored group names as text
# ("G1".."G7"). Switch to plain numeric group numbers (1..7) so a group can be
# referenced/matched numerically when assigning multiple professors per module.
$wsL3 = $wb.Worksheets.Item("L3")
$wsL3.Range("A2").Value = 1
$wsL3.Range("A3").Value = 2
$wsL3.Range("A4").Value = 3
$wsL3.Range("A5").Value = 4
$wsL3.Range("A6").Value = 5
$wsL3.Range("A7").Value = 6
$wsL3.Range("A8").Value = 7

# Update the selection left on the "Professors" sheet.
$wsProf = $wb.Worksheets.Item("Professors")
$wsProf.Range("C7").Select()

# Leave the "L3" sheet as the active tab with its own selection, matching the
# state the workbook was saved in.
$wsL3.Range("A9").Select()
$wsL3.Activate()
